$d = $word.ActiveDocument

# --- Step 1: remove the whole "Meta description: ..." paragraph ---
$metaFind = $d.Content
$metaFind.Find.ClearFormatting()
$metaOk = $metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$metaPara = $metaFind.Duplicate
$metaPara.Expand(4) | Out-Null
$metaPara.Delete()

# --- Step 2: the last paragraph (the italic "Can I play ... mobile device?"
#     blurb) is replaced by two new paragraphs: a bold heading-style line
#     repeating the page title, followed by the italic meta-description
#     blurb that used to live at the top of the document. ---
$pCount = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($pCount)
$pLast.Range.Delete()

$pCount2 = $d.Paragraphs.Count
$pLast2 = $d.Paragraphs.Item($pCount2)
$endPos = $pLast2.Range.End
$insPoint = $d.Range($endPos, $endPos)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Bucks Buffalo Gigablox for Free – Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Explore the Western desert with Big Bucks Buffalo Gigablox – read our review, then play for free here! Gigablox symbols and Scatter Respins add excitement.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insPoint.InsertXML($xml) | Out-Null
